# "Full and partial matching routine"
#
# Rows 56-62 of the WORLDLS sheet had the French ("FR", column E) and
# Catalan ("CA", column D) weekday translations swapped - column D held the
# French word and column E held the Catalan word, the opposite of every
# other row in the table. This fixes the D/E values for those seven rows
# (Monday..Sunday) by swapping them back into the correct columns, and
# restores the sheet's scroll/selection state left over from reviewing the
# fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the mismatched Catalan (D) / French (E) pair for each weekday row.
# (NB: read via .Value2 - the plain .Value getter misbehaves in this host.)
56..62 | ForEach-Object {
    $row = $_
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)

    $dValue = $dCell.Value2
    $eValue = $eCell.Value2

    $dCell.Value2 = $eValue
    $eCell.Value2 = $dValue
}

# Leave the view the way it was after making/reviewing the fix.
$ws.Range("D67").Select()
